$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C ("denomination") between B (Article name) and the old C (taille/prix)
$ws.Columns.Item(3).Insert()

# Match column C width/format roughly to column B (closest achievable via ColumnWidth)
$ws.Columns.Item(3).ColumnWidth = 22.83

# Header
$ws.Range("C1").Value = "denomination"

# Denomination values per row
$ws.Cells.Item(2, 3).Value = "1/8"
$ws.Cells.Item(3, 3).Value = "1/4"
$ws.Cells.Item(4, 3).Value = "3/8"
$ws.Cells.Item(5, 3).Value = "1/2"
$ws.Cells.Item(6, 3).Value = "3/4"
$ws.Cells.Item(7, 3).Value = "4/4"
$ws.Cells.Item(8, 3).Value = "5/4"
$ws.Cells.Item(9, 3).Value = "6/4"
$ws.Cells.Item(10, 3).Value = "1/8"
$ws.Cells.Item(11, 3).Value = "1/4"
$ws.Cells.Item(12, 3).Value = "3/8"
$ws.Cells.Item(13, 3).Value = "1/2"
$ws.Cells.Item(14, 3).Value = "3/4"
$ws.Cells.Item(15, 3).Value = "4/4"
$ws.Cells.Item(16, 3).Value = "5/4"
$ws.Cells.Item(17, 3).Value = "6/4"
$ws.Cells.Item(18, 3).Value = "3/8"
$ws.Cells.Item(19, 3).Value = "1/2"
$ws.Cells.Item(20, 3).Value = "3/4"
$ws.Cells.Item(21, 3).Value = "4/4"
$ws.Cells.Item(22, 3).Value = "5/4"
$ws.Cells.Item(23, 3).Value = "6/4"
$ws.Cells.Item(24, 3).Value = "3/8"
$ws.Cells.Item(25, 3).Value = "1/2"
$ws.Cells.Item(26, 3).Value = "3/4"
$ws.Cells.Item(27, 3).Value = "4/4"
$ws.Cells.Item(28, 3).Value = "5/4"
$ws.Cells.Item(29, 3).Value = "6/4"
$ws.Cells.Item(30, 3).Value = "1/2"
$ws.Cells.Item(31, 3).Value = "3/4"
$ws.Cells.Item(32, 3).Value = "1/4"
$ws.Cells.Item(33, 3).Value = "3/8"
$ws.Cells.Item(34, 3).Value = "1/2"
$ws.Cells.Item(35, 3).Value = "3/4 "
$ws.Cells.Item(36, 3).Value = "4/4"
$ws.Cells.Item(37, 3).Value = "5/4"
$ws.Cells.Item(38, 3).Value = "6/4"
$ws.Cells.Item(39, 3).Value = "1/2"
$ws.Cells.Item(40, 3).Value = "3/4"
$ws.Cells.Item(41, 3).Value = "1/8"
$ws.Cells.Item(42, 3).Value = "1/4"
$ws.Cells.Item(43, 3).Value = "3/8"
$ws.Cells.Item(44, 3).Value = "1/2"
$ws.Cells.Item(45, 3).Value = "3/4"
$ws.Cells.Item(46, 3).Value = "4/4"
$ws.Cells.Item(47, 3).Value = "5/4"
$ws.Cells.Item(48, 3).Value = "6/4"
$ws.Cells.Item(49, 3).Value = "1/4-1/8"
$ws.Cells.Item(50, 3).Value = "3/8-1/4"
$ws.Cells.Item(51, 3).Value = "3/8-1/8"
$ws.Cells.Item(52, 3).Value = " 1/2-3/8"
$ws.Cells.Item(53, 3).Value = "3/4-1/2"
$ws.Cells.Item(54, 3).Value = "4/4-3/4"
$ws.Cells.Item(55, 3).Value = "5/4-4/4"
$ws.Cells.Item(56, 3).Value = "6/4-5/4"
$ws.Cells.Item(57, 3).Value = "5/4-3/4"
$ws.Cells.Item(59, 3).Value = "1/4M-1/8F"
$ws.Cells.Item(60, 3).Value = "3/8M-1/8F"
$ws.Cells.Item(61, 3).Value = "3/8M-1/4F"
$ws.Cells.Item(62, 3).Value = "1/2M-3/8F"
$ws.Cells.Item(63, 3).Value = "3/4M-1/2F"
$ws.Cells.Item(64, 3).Value = "4/4M-1/2F"
$ws.Cells.Item(65, 3).Value = "4/4M-3/4F"
$ws.Cells.Item(66, 3).Value = "5/4M-4/4F"
$ws.Cells.Item(67, 3).Value = "6/4M-5/4F"
$ws.Cells.Item(69, 3).Value = "1/4F-1M"
$ws.Cells.Item(70, 3).Value = "3/8F-1M"
$ws.Cells.Item(71, 3).Value = "3/8F-1M"
$ws.Cells.Item(72, 3).Value = "1/2F-1M"
$ws.Cells.Item(73, 3).Value = "1/2F-3M"
$ws.Cells.Item(74, 3).Value = "3/4F-1M"
$ws.Cells.Item(75, 3).Value = "4/4F-1M"
$ws.Cells.Item(76, 3).Value = "4/4F-3M"
$ws.Cells.Item(77, 3).Value = "5/4F-4M"
$ws.Cells.Item(78, 3).Value = "6/4F-5M"
$ws.Cells.Item(79, 3).Value = "1/4-1/8"
$ws.Cells.Item(80, 3).Value = "3/8-1/4"
$ws.Cells.Item(81, 3).Value = "1/2-3/8"
$ws.Cells.Item(82, 3).Value = "3/4-1/2"
$ws.Cells.Item(83, 3).Value = "4/4-3/4"
$ws.Cells.Item(84, 3).Value = "4/4-1/2"
$ws.Cells.Item(85, 3).Value = "5/4-4/4"
$ws.Cells.Item(86, 3).Value = "6/4-4/4"
$ws.Cells.Item(87, 3).Value = "6/4-5/4"
$ws.Cells.Item(89, 3).Value = "1/8"
$ws.Cells.Item(90, 3).Value = "1/4"
$ws.Cells.Item(91, 3).Value = "3/8"
$ws.Cells.Item(92, 3).Value = "1/2"
$ws.Cells.Item(93, 3).Value = "3/4-10"
$ws.Cells.Item(94, 3).Value = "4/4"
$ws.Cells.Item(95, 3).Value = "5/4"
$ws.Cells.Item(96, 3).Value = "6/4"
$ws.Cells.Item(97, 3).Value = "1/2-2"
$ws.Cells.Item(98, 3).Value = "1/4"
$ws.Cells.Item(99, 3).Value = "3/8"
$ws.Cells.Item(100, 3).Value = "1/2"
$ws.Cells.Item(101, 3).Value = "3/4"
$ws.Cells.Item(102, 3).Value = "4/4"
$ws.Cells.Item(103, 3).Value = "5/4"
$ws.Cells.Item(104, 3).Value = "6/4"
$ws.Cells.Item(105, 3).Value = "3/8"
$ws.Cells.Item(106, 3).Value = "1/2"
$ws.Cells.Item(107, 3).Value = "3/4"
$ws.Cells.Item(108, 3).Value = "4/4"
$ws.Cells.Item(109, 3).Value = "5/4"
$ws.Cells.Item(110, 3).Value = "6/4"
$ws.Cells.Item(111, 3).Value = "1/2-12"
$ws.Cells.Item(112, 3).Value = "1/2-15"
$ws.Cells.Item(113, 3).Value = "1/2-1/2"
$ws.Cells.Item(114, 3).Value = "1/2-1/2-1EF"
$ws.Cells.Item(115, 3).Value = "15-1/2-15"
$ws.Cells.Item(116, 3).Value = "3/4-22"
$ws.Cells.Item(117, 3).Value = "3/4"
$ws.Cells.Item(118, 3).Value = "1/2-10"
$ws.Cells.Item(119, 3).Value = "1/2-20"
$ws.Cells.Item(120, 3).Value = "1/2-30"
$ws.Cells.Item(121, 3).Value = "1/2-50"
$ws.Cells.Item(122, 3).Value = "1/2-100"
$ws.Cells.Item(123, 3).Value = "3/4-15"
$ws.Cells.Item(124, 3).Value = "3/4-20"
$ws.Cells.Item(125, 3).Value = "3/4-30"
$ws.Cells.Item(126, 3).Value = "3/4-50"
$ws.Cells.Item(127, 3).Value = "3/4-100"

# Restore selection to match the authored state
$ws.Range("C1:C127").Select()
